# Insert a new weekly data row at row 11 (pushing the existing rows 11-68
# down to 12-69, matching the "Fruta / hortaliza, semanal" update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new weekly record.
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "Vega Monumental Concepción"
$ws.Range("C11").Value = "Bíobío"
$ws.Range("D11").Value = 44600
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 100112001
$ws.Range("G11").Value = "Berenjena"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 7500
$ws.Range("L11").Value = 8000
$ws.Range("M11").Value = 7760
$ws.Range("N11").Value = "$/caja 60 unidades"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 129
$ws.Range("Q11").Value = 60
$ws.Range("R11").Value = "Hortaliza"
